# "4.1V Full charge 수정" - add Source Path info row and two new firmware
# revision history rows (V2.10_1.0sec / V2.10_0.5sec) to the revision table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New row 2: Source Path label + value -------------------------------
$ws.Range("B2").Value = "Source Path"
$ws.Range("C2").Value = "C:\WORK\Project\3_Source\Test Sourse\V2.0_Source_20190403\plasma_pt"

# B2 keeps the column's default (centered) look; C2 holds a long path so it
# reads better left aligned.
$ws.Range("C2").HorizontalAlignment = -4131  # xlLeft
$ws.Range("C2").VerticalAlignment = -4108    # xlCenter

# --- Row 10: new firmware revision entry --------------------------------
$revDate = Get-Date -Year 2019 -Month 5 -Day 20 -Hour 0 -Minute 0 -Second 0
$ws.Range("B10").Value = $revDate
$ws.Range("C10").Value = "V2.10_1.0sec"
$ws.Range("D10").Value = "GAS_EN 후 1.0sec 후 plasma On, 4.1V 완충(ADC238)"

# --- Row 11: new firmware revision entry --------------------------------
$ws.Range("B11").Value = $revDate
$ws.Range("C11").Value = "V2.10_0.5sec"
$ws.Range("D11").Value = "GAS_EN 후 0.5sec 후 plasma On, 4.1V 완충(ADC238)"

# Match the author's last active selection when the file was saved.
$ws.Range("D13").Select() | Out-Null
